$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("G1")

$scratch.Formula = "=""35.346.59"""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  +0.33%  "
$scratch.Formula = "=""1.901.15"""
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  +2.05%  "
$scratch.Formula = "=""245.78"""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +2.64%  "
$scratch.Formula = "=""0.663"""
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +6.29%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("E9").Value = "  +5.56%  "
$scratch.Formula = "=""52.96"""
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +12.83%  "
$ws.Range("E11").Value = "  +3.41%  "
$scratch.Formula = "=""0.0995"""
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +0.48%  "
$scratch.Formula = "=""2.176.95"""
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +2.13%  "
$scratch.Formula = "=""12.11"""
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +5.33%  "
$ws.Range("E15").Value = "  +2.48%  "
$scratch.Formula = "=""1.900.98"""
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +2.43%  "
$scratch.Formula = "=""35.323.53"""
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +0.34%  "
$scratch.Formula = "=""72.22"""
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +3.21%  "
$scratch.Formula = "=""0.0₃0818"""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +2.73%  "
$scratch.Formula = "=""240.65"""
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("E24").Value = "  -0.35%  "
$scratch.Formula = "=""2.30"""
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("E26").Value = "  +22.67%  "
$scratch.Formula = "=""170.09"""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +3.17%  "
$scratch.Formula = "=""18.33"""
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("E31").Value = "  +2.68%  "
$scratch.Formula = "=""0.0564"""
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +0.16%  "
$scratch.Formula = "=""0.931"""
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +13.49%  "
$scratch.Formula = "=""4.08"""
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("E37").Value = "  -2.12%  "
$scratch.Formula = "=""1.33"""
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("E39").Value = "  -0.93%  "
$scratch.Formula = "=""0.0208"""
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +2.25%  "
$scratch.Formula = "=""16.06"""
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +5.73%  "
$scratch.Formula = "=""0.0626"""
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +4.06%  "
$scratch.Formula = "=""89.51"""
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.04%  "
$scratch.Formula = "=""1.335.88"""
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.85%  "
$scratch.Formula = "=""2.36"""
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +1.17%  "
$scratch.Formula = "=""48.00"""
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +38.11%  "
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  -1.39%  "
$scratch.Formula = "=""11.83"""
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -4.82%  "
$scratch.Formula = "=""2.085.88"""
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +1.97%  "

$scratch.Clear()
$excel.CutCopyMode = $false
